# Allow user to change port of dash app
#
# - "settings" sheet, cell B4 ("case"): "D" -> "BPV"
# - "settings" sheet: add a new "port" row (param/value/info) with value 8050

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Update the existing "case" row value from "D" to "BPV"
$ws.Range("B4").Value = "BPV"

# Append a new row describing the dash app's port
$ws.Range("A6").Value = "port"
$ws.Range("B6").Value = 8050
$ws.Range("C6").Value = "port number on which the dash app displays"
